$d = $word.ActiveDocument

# The first two paragraphs:
#   1) "House of Hospitality" (italic) + "," + " " + <line break> +
#      "Chapter Five ========================="
#   2) "By Dorothy Day" (bold)
# are being collapsed into a single, unformatted pandoc-style title-block
# paragraph: "% Dorothy Day"

$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)

# Remove both paragraphs completely (text + paragraph marks).
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$rng.Delete()

# Insert a brand-new, unformatted paragraph in their place and give it the
# new text. Using InsertParagraphBefore + a fresh Text assignment (rather
# than typing into the old, now-deleted runs) ensures no leftover
# italic/bold character formatting is carried over.
$insertionPoint = $d.Range($d.Paragraphs(1).Range.Start, $d.Paragraphs(1).Range.Start)
$insertionPoint.InsertParagraphBefore()

$newPara = $d.Paragraphs(1).Range
$newParaText = $d.Range($newPara.Start, $newPara.End - 1)
$newParaText.Text = "% Dorothy Day"
